$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.987.50"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "2.047.87"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.33"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("E6").Value = "  +2.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.11"
$ws.Range("E7").Value = "  +7.17%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.380"
$ws.Range("E9").Value = "  +2.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0778"
$ws.Range("E10").Value = "  +3.69%  "

$ws.Range("E11").Value = "  +2.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.77"
$ws.Range("E12").Value = "  +6.36%  "

$ws.Range("D13").Value = "2.349.51"
$ws.Range("E13").Value = "  +0.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.800"
$ws.Range("E14").Value = "  -0.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.56"
$ws.Range("E15").Value = "  +8.49%  "

$ws.Range("D16").Value = "2.049.55"
$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("D17").Value = "37.039.30"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.59"
$ws.Range("E18").Value = "  +18.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.47"
$ws.Range("E19").Value = "  +4.44%  "

$ws.Range("D20").Value = "0.0₃0900"
$ws.Range("E20").Value = "  +1.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.34"
$ws.Range("E21").Value = "  +3.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.60"
$ws.Range("E22").Value = "  +0.34%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("E24").Value = "  -0.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  +11.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.60"
$ws.Range("E26").Value = "  -0.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.09"
$ws.Range("E27").Value = "  +2.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.72"
$ws.Range("E28").Value = "  -0.18%  "

$ws.Range("E29").Value = "  +2.03%  "

$ws.Range("E30").Value = "  +8.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.68"
$ws.Range("E31").Value = "  +4.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0611"
$ws.Range("E32").Value = "  +0.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.46"
$ws.Range("E33").Value = "  +4.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0881"
$ws.Range("E34").Value = "  +1.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.21"
$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("E37").Value = "  -2.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.107"
$ws.Range("E38").Value = "  +5.16%  "

$ws.Range("E39").Value = "  +0.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.18"
$ws.Range("E40").Value = "  +14.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.92"
$ws.Range("E41").Value = "  +25.08%  "

$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("E43").Value = "  -3.75%  "

$ws.Range("E44").Value = "  -0.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.48"
$ws.Range("E45").Value = "  +1.01%  "

$ws.Range("E46").Value = "  +5.53%  "

$ws.Range("D47").Value = "1.279.26"
$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("E48").Value = "  -1.25%  "

$ws.Range("D49").Value = "2.236.77"
$ws.Range("E49").Value = "  +0.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.66"
$ws.Range("E50").Value = "  -0.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.60"
$ws.Range("E51").Value = "  -9.01%  "
